$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tab10")

# Fix mojibake in the Regional Economic Communities footnote (A103):
# replace corrupted "Pa>ses"/"L>ngua"/"Com>n" sequences with the correct
# accented Portuguese/Spanish characters.
$ws.Range("A103").Value = 'Regional Economic Communities:CEN-SAD = "Community of Sahel-Saharan States";COMESA = "Common Market for Eastern and Southern Africa";EAC = "East African Community";ECCAS = "Economic Community of Central African States";ECOWAS = "Economic Community of West African States";IGAD = "Intergovernmental Authority on Development";SADC = "Southern African Development Community";UMA = "Arab Maghreb Union";PALOP = "Países Africanos de Língua Oficial Portuguesa";ASEAN = "Association of Southeast Asian Nations";MERCOSUR = "Mercado Común del Sur".EU27 = "European Union (27 members)".OECD = "Organisation for Economic Co-operation and Development".'

# Update a handful of statistical data values that were revised.
$ws.Range("I68").Value = 45.645454545454598

$ws.Range("C97").Value = 0.52607142857142997
$ws.Range("D97").Value = 0.88707692307691999
$ws.Range("E97").Value = 0.56604545454544997
$ws.Range("F97").Value = 39.829166666666701
$ws.Range("G97").Value = 52.839285714285701
$ws.Range("H97").Value = 31.02
$ws.Range("I97").Value = 41.588000000000001
$ws.Range("J97").Value = 34.570370370370398

$ws.Range("C98").Value = 0.63606666666667
$ws.Range("D98").Value = 0.91858333333332998
$ws.Range("E98").Value = 0.48275000000000001
$ws.Range("F98").Value = 36.036363636363603
$ws.Range("G98").Value = 41.4
$ws.Range("H98").Value = 34.613333333333301
$ws.Range("I98").Value = 35.200000000000003
$ws.Range("J98").Value = 29.072727272727299
